$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# --- Update car-time values around the Delft area ---
# (G column = Delft; row 7 = Delft; symmetric matrix, the TRANSPOSE()
#  array formulas on the other half recompute automatically)
$ws.Range("G8").Value  = 24   # Den Haag    -> Delft   (was 20)
$ws.Range("N10").Value = 28   # Gouda       -> Rotterdam (was 25)
$ws.Range("G14").Value = 20   # Rotterdam   -> Delft   (was 19)
$ws.Range("J14").Value = 28   # Rotterdam   -> Haarlem (was 25)
$ws.Range("G19").Value = 20   # Vlaardingen -> Delft   (was 18)
$ws.Range("G20").Value = 23   # Westland    -> Delft   (was 19)
$ws.Range("S20").Value = 16   # Westland    -> Vlaardingen (was 19)

# --- Toggle A15 interior pattern so the cached "apply fill" flag drops
#     (A15's fill stays "none", matching its original appearance) ---
$a15 = $ws.Range("A15")
$a15.Interior.Pattern = 17      # transient pattern
$a15.Interior.Pattern = -4142   # xlNone  (back to original look)

# --- Header formatting: right-align the "Delft" header cell (G1) ---
$ws.Range("G1").HorizontalAlignment = -4152   # xlRight

# --- Move the active selection to G20 ---
$ws.Range("G20").Select()
